$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: A15 = "06-10-2025", B15 = gold-price sentence for 06-10-2025 ---

# A15 needs to be a *text* cell (not auto-converted to a date serial). Stage
# the literal text in a scratch cell forced to Text format, then bring only
# the resulting string value into A15 (whose style we first clone from A14)
# via a values-only paste so A15's style/number-format never changes.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "06-10-2025"

$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4104)

$ws.Range("ZZ1").Copy()
$ws.Range("A15").PasteSpecial(-4163)

$ws.Range("ZZ1").Clear()

# B15 holds free-form text that never looks like a date/number, so a direct
# assignment after cloning B14's style is safe.
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4104)
$ws.Range("B15").Value = "The price of gold in India today is ₹12,077 per gram for 24 karat gold, ₹11,070 per gram for 22 karat gold and ₹9,058 per gram for 18 karat gold (also called 999 gold)."
